# Actualización automática 2025-10-16 08:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M31").Value = 2296.77
$ws1.Range("M36").Value = "7 de 34"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F31").Value = 2296.77
$ws2.Range("F36").Value = 16505.9

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 15092.49
$ws3.Range("E12").Value = 6608.780000000001
$ws3.Range("F12").Value = 0.6954657492395606

$ws3.Range("D14").Value = 16505.9
$ws3.Range("E14").Value = 20079.66723718182
$ws3.Range("F14").Value = 0.4511587832708274
